$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.872.29"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "1.622.98"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").Value = "'0.996"
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").Value = "'214.21"
$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").Value = "'29.75"
$ws.Range("E8").Value = "  +10.97%  "
$ws.Range("D9").Value = "'0.258"
$ws.Range("E9").Value = "  +2.80%  "
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.856.63"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "1.623.15"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").Value = "'0.568"
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("E15").Value = "  +4.80%  "
$ws.Range("D16").Value = "29.930.56"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").Value = "'8.79"
$ws.Range("E17").Value = "  +15.75%  "
$ws.Range("D18").Value = "'64.59"
$ws.Range("E18").Value = "  +1.85%  "
$ws.Range("D19").Value = "'243.97"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +3.17%  "
$ws.Range("D23").Value = "'9.61"
$ws.Range("E23").Value = "  +4.19%  "
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").Value = "'157.04"
$ws.Range("E26").Value = "  +2.52%  "
$ws.Range("D27").Value = "'0.110"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  -0.32%  "
$ws.Range("E30").Value = "  +3.39%  "
$ws.Range("E31").Value = "  +5.23%  "
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("D34").Value = "1.426.77"
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("E35").Value = "  +6.89%  "
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +2.94%  "
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").Value = "'0.0507"
$ws.Range("E41").Value = "  +3.22%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").Value = "'0.832"
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").Value = "'53.95"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'69.16"
$ws.Range("E45").Value = "  +4.76%  "
$ws.Range("E46").Value = "  +18.34%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'5.42"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "1.765.00"
$ws.Range("E49").Value = "  +1.03%  "
$ws.Range("D50").Value = "'88.40"
$ws.Range("E50").Value = "  +2.04%  "
$ws.Range("E51").Value = "  +2.49%  "
